$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are treated as text so that numeric-looking values
# (e.g. "18.00", "5.00") are not auto-converted to numbers and keep their
# exact formatting/trailing zeros, matching the original inline-string cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.903.04"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.307.77"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.57"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.29"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.299.50"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.626"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.24"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.07%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.835.10"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.120"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.93%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.00"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.320.75"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "63.851.48"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.977"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "457.87"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.00"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +8.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.03"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.53%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"

$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.79"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.79%  "

$ws.Range("B26").Value = "Litecoin"

$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "86.26"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.42%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.82%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.49"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.21%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "60.95"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "560.88"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.34%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.04%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.89%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.99"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.00%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0721"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.033.43"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.28%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.15%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.18"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.35%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.01%  "

$ws.Range("B48").Value = "FirstDigitalUSD"

$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.18%  "

$ws.Range("B49").Value = "Monero"

$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.74"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.50"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.06"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.22%  "
